$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 456.54544
$ws.Range("J33").Value = 477.4
$ws.Range("L33").Value = 477.4
$ws.Range("N33").Value = -935.4
$ws.Range("H98").Value = 2349.476
$ws.Range("I98").Value = 1822.875
$ws.Range("J98").Value = 4034.6
$ws.Range("K98").Value = 1822.875
$ws.Range("L98").Value = 4034.6
$ws.Range("M98").Value = -324.875
$ws.Range("N98").Value = -7030.6
$ws.Range("H112").Value = 1296.7894
$ws.Range("J112").Value = 1297.027
$ws.Range("L112").Value = 3891.081
$ws.Range("N112").Value = -6107.081
$ws.Range("H113").Value = 334791.34
$ws.Range("J113").Value = 2369
$ws.Range("L113").Value = 2369
$ws.Range("N113").Value = -8877
$ws.Range("H115").Value = 761.2857
$ws.Range("I115").Value = 761.2857
$ws.Range("K115").Value = 2283.8571
$ws.Range("M115").Value = -716.8571000000002
$ws.Range("H122").Value = 2349.476
$ws.Range("I122").Value = 1822.875
$ws.Range("J122").Value = 4034.6
$ws.Range("K122").Value = 5468.625
$ws.Range("L122").Value = 12103.8
$ws.Range("M122").Value = -3018.625
$ws.Range("N122").Value = -17003.8
$ws.Range("H132").Value = 1833.3572
$ws.Range("I132").Value = 1510.7894
$ws.Range("K132").Value = 4532.3682
$ws.Range("M132").Value = -2002.3682
$ws.Range("H135").Value = 77658.08
$ws.Range("I135").Value = 796.5
$ws.Range("K135").Value = 7168.5
$ws.Range("M135").Value = -4633.5
$ws.Range("H137").Value = 1911.1364
$ws.Range("I137").Value = 1969.7222
$ws.Range("J137").Value = 1647.5
$ws.Range("K137").Value = 5909.1666
$ws.Range("L137").Value = 4942.5
$ws.Range("M137").Value = -3359.1666
$ws.Range("N137").Value = -10042.5
$ws.Range("H138").Value = 1710.5172
$ws.Range("I138").Value = 1075.381
$ws.Range("J138").Value = 3377.75
$ws.Range("K138").Value = 3226.143
$ws.Range("L138").Value = 10133.25
$ws.Range("M138").Value = 1913.857
$ws.Range("N138").Value = -20413.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23763.434
$ws.Range("I32").Value = 4989.673
$ws.Range("K32").Value = 4989.673
$ws.Range("M32").Value = -4702.673
$ws.Range("H45").Value = 5360.8823
$ws.Range("I45").Value = 6751.55
$ws.Range("K45").Value = 6751.55
$ws.Range("M45").Value = -6374.55
$ws.Range("H132").Value = 1191.4529
$ws.Range("I132").Value = 1158.7551
$ws.Range("K132").Value = 3476.2653
$ws.Range("M132").Value = -946.2653

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5067.5
$ws.Range("I20").Value = 4398.35
$ws.Range("J20").Value = 6023.4287
$ws.Range("K20").Value = 4398.35
$ws.Range("L20").Value = 6023.4287
$ws.Range("M20").Value = -4151.35
$ws.Range("N20").Value = -6517.4287
$ws.Range("H99").Value = 3972.0667
$ws.Range("I99").Value = 4048.4167
$ws.Range("K99").Value = 4048.4167
$ws.Range("M99").Value = -2550.4167
$ws.Range("H134").Value = 905.7273
$ws.Range("I134").Value = 941.3
$ws.Range("K134").Value = 2823.9
$ws.Range("M134").Value = -288.8999999999996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13014.973
$ws.Range("J31").Value = 22252.5
$ws.Range("L31").Value = 22252.5
$ws.Range("N31").Value = -22842.5
$ws.Range("H34").Value = 13014.973
$ws.Range("J34").Value = 22252.5
$ws.Range("L34").Value = 22252.5
$ws.Range("N34").Value = -22656.5
$ws.Range("H99").Value = 17465.375
$ws.Range("I99").Value = 27007
$ws.Range("K99").Value = 27007
$ws.Range("M99").Value = -25509
$ws.Range("H107").Value = 1383.6904
$ws.Range("I107").Value = 1294.7916
$ws.Range("J107").Value = 1502.2222
$ws.Range("K107").Value = 1294.7916
$ws.Range("L107").Value = 1502.2222
$ws.Range("M107").Value = 625.2084
$ws.Range("N107").Value = -5342.2222
$ws.Range("H126").Value = 17465.375
$ws.Range("I126").Value = 27007
$ws.Range("K126").Value = 81021
$ws.Range("M126").Value = -78551
$ws.Range("H132").Value = 2133.9285
$ws.Range("I132").Value = 1922.04
$ws.Range("K132").Value = 5766.12
$ws.Range("M132").Value = -3236.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 552
$ws.Range("I2").Value = 818.25
$ws.Range("K2").Value = 4909.5
$ws.Range("M2").Value = -4796.5
$ws.Range("H117").Value = 5266.778
$ws.Range("I117").Value = 2249
$ws.Range("J117").Value = 6427.4614
$ws.Range("K117").Value = 6747
$ws.Range("L117").Value = 19282.3842
$ws.Range("M117").Value = -3305
$ws.Range("N117").Value = -26166.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H97").Value = 14232.0205
$ws.Range("I97").Value = 20042.588
$ws.Range("K97").Value = 20042.588
$ws.Range("M97").Value = -19546.588
$ws.Range("H102").Value = 1496.5
$ws.Range("I102").Value = 1519.0435
$ws.Range("K102").Value = 1519.0435
$ws.Range("M102").Value = 102.9565
$ws.Range("H132").Value = 3957.8635
$ws.Range("I132").Value = 3754.625
$ws.Range("J132").Value = 4499.8335
$ws.Range("K132").Value = 11263.875
$ws.Range("L132").Value = 13499.5005
$ws.Range("M132").Value = -8733.875
$ws.Range("N132").Value = -18559.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9679.541999999999
$ws.Range("I7").Value = 16292.2
$ws.Range("J7").Value = 4956.2144
$ws.Range("K7").Value = 16292.2
$ws.Range("L7").Value = 4956.2144
$ws.Range("M7").Value = -16180.2
$ws.Range("N7").Value = -5180.2144
$ws.Range("H22").Value = 2051.5557
$ws.Range("I22").Value = 571.3333
$ws.Range("J22").Value = 2791.6667
$ws.Range("K22").Value = 571.3333
$ws.Range("L22").Value = 2791.6667
$ws.Range("M22").Value = -276.3333
$ws.Range("N22").Value = -3381.6667
$ws.Range("H27").Value = 2051.5557
$ws.Range("I27").Value = 571.3333
$ws.Range("J27").Value = 2791.6667
$ws.Range("K27").Value = 571.3333
$ws.Range("L27").Value = 2791.6667
$ws.Range("M27").Value = -464.3333
$ws.Range("N27").Value = -3005.6667
$ws.Range("H40").Value = 3607.95
$ws.Range("I40").Value = 2885.9375
$ws.Range("K40").Value = 2885.9375
$ws.Range("M40").Value = -2749.9375
$ws.Range("H122").Value = 8194.9
$ws.Range("I122").Value = 10349.077
$ws.Range("J122").Value = 4194.2856
$ws.Range("K122").Value = 31047.231
$ws.Range("L122").Value = 12582.8568
$ws.Range("M122").Value = -28597.231
$ws.Range("N122").Value = -17482.8568
$ws.Range("H126").Value = 9679.541999999999
$ws.Range("I126").Value = 16292.2
$ws.Range("J126").Value = 4956.2144
$ws.Range("K126").Value = 48876.60000000001
$ws.Range("L126").Value = 14868.6432
$ws.Range("M126").Value = -46406.60000000001
$ws.Range("N126").Value = -19808.6432
$ws.Range("H132").Value = 2856.9
$ws.Range("I132").Value = 2171.32
$ws.Range("J132").Value = 3999.5334
$ws.Range("K132").Value = 6513.960000000001
$ws.Range("L132").Value = 11998.6002
$ws.Range("M132").Value = -3983.960000000001
$ws.Range("N132").Value = -17058.6002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2441.7058
$ws.Range("I126").Value = 2138.3845
$ws.Range("K126").Value = 6415.1535
$ws.Range("M126").Value = -3945.1535
$ws.Range("H132").Value = 1132956.5
$ws.Range("I132").Value = 1321361.1
$ws.Range("J132").Value = 2528.8
$ws.Range("K132").Value = 3964083.3
$ws.Range("L132").Value = 7586.400000000001
$ws.Range("M132").Value = -3961553.3
$ws.Range("N132").Value = -12646.4
$ws.Range("H136").Value = 632.8889
$ws.Range("I136").Value = 632.8889
$ws.Range("K136").Value = 1898.6667
$ws.Range("M136").Value = 651.3332999999998
